$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 28.5

# Row 3
$ws.Range("B3").Value = 3033
$ws.Range("C3").Value = -5850
$ws.Range("D3").Value = 31

# Row 4 (was text "233.50", now plain number)
$ws.Range("B4").Value = 3035
$ws.Range("C4").Value = -5920
$ws.Range("D4").Value = 191

# Row 5 (was plain number 250, now text "256.50")
$ws.Range("B5").Value = 3141
$ws.Range("C5").Value = -5899.5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.50"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("B6").Value = 3189
$ws.Range("C6").Value = -6345
$ws.Range("D6").Value = 243.5

# Row 7 (was text "48.10", now plain number)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 3504
$ws.Range("C7").Value = -5934
$ws.Range("D7").Value = 130

# Row 8 (stays text, content changes)
$ws.Range("A8").Value = 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "257.00"
$ws.Range("D8").Style = "Normal"
